# Insert a new weekly record at row 545, shifting the existing rows
# 545-573 down to 546-574 (same pattern as every prior weekly update in
# this sheet: newest observation inserted just above the rest of the
# "Femacal de La Calera" block, everything below slides down by one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 545:573 down to 546:574, leaving a blank row 545 to fill in.
$ws.Rows(545).Insert()

# Populate the new row 545 with the new weekly observation.
$ws.Cells.Item(545, 1).Value = 3
$ws.Cells.Item(545, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(545, 3).Value = "Coquimbo"
$ws.Cells.Item(545, 4).Value = 45041
$ws.Cells.Item(545, 5).Value = 5
$ws.Cells.Item(545, 6).Value = 100112031
$ws.Cells.Item(545, 7).Value = "Poroto verde"
$ws.Cells.Item(545, 8).Value = "Magnum"
$ws.Cells.Item(545, 9).Value = "Primera"
$ws.Cells.Item(545, 10).Value = 82
$ws.Cells.Item(545, 11).Value = 29000
$ws.Cells.Item(545, 12).Value = 30000
$ws.Cells.Item(545, 13).Value = 29512
$ws.Cells.Item(545, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(545, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(545, 16).Value = 1180
$ws.Cells.Item(545, 17).Value = 25
$ws.Cells.Item(545, 18).Value = "Hortaliza"
